$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (AC1) onto the new header cells
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill in the team record (Wins/Losses/Ties) for every data row (2-52)
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 93
    $ws.Cells.Item($r, 31).Value = 69
    $ws.Cells.Item($r, 32).Value = 0
}
